$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 5 (pushes the old row5.."1",2000.. etc down by one)
$ws.Rows("5:5").Insert(-4121)

# The inserted row 5 becomes the new "99999" row (like A2:A4, A is blank/numeric)
$ws.Cells.Item(5, 1).Value = 0
$ws.Cells.Item(5, 2).Value = 99999

# Old row4 (still row 4) changes its B value from 99999 to 99993
$ws.Cells.Item(4, 2).Value = 99993
